$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-8 (ME, MAE, MSE, RMSE, SE) -------------------

$ws.Range("B2").Value = 0.223982636080522
$ws.Range("C2").Value = 0.7951022744350149
$ws.Range("D2").Value = 1.433395224409647
$ws.Range("E2").Value = 1.197244847309708
$ws.Range("F2").Value = 1.220503617491155

$ws.Range("B3").Value = 0.03460922587986636
$ws.Range("C3").Value = 0.5169290182150836
$ws.Range("D3").Value = 0.4531020548674881
$ws.Range("E3").Value = 0.6731285574594856
$ws.Range("F3").Value = 0.699687747770943

$ws.Range("B4").Value = 0.05542554409107694
$ws.Range("C4").Value = 0.5672532340729006
$ws.Range("D4").Value = 0.5863441778978949
$ws.Range("E4").Value = 0.765731139433349
$ws.Range("F4").Value = 0.7976822244661111

$ws.Range("B5").Value = 0.09157044451040193
$ws.Range("C5").Value = 0.6935010056718082
$ws.Range("D5").Value = 0.734265071824128
$ws.Range("E5").Value = 0.8568926839599741
$ws.Range("F5").Value = 0.8935703207179192

$ws.Range("B6").Value = 0.1506695717468531
$ws.Range("C6").Value = 0.7208337799200664
$ws.Range("D6").Value = 0.8151397471296369
$ws.Range("E6").Value = 0.9028508997224497
$ws.Range("F6").Value = 0.9383427632914298

$ws.Range("B7").Value = 0.1498722168069911
$ws.Range("C7").Value = 0.8124485921203141
$ws.Range("D7").Value = 1.075051427048039
$ws.Range("E7").Value = 1.036846867694569
$ws.Range("F7").Value = 1.088192751256446
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = 0.09843046095639814
$ws.Range("C8").Value = 0.9381895697043964
$ws.Range("D8").Value = 1.212929050766225
$ws.Range("E8").Value = 1.101330581962666
$ws.Range("F8").Value = 1.201619155201249
$ws.Range("G8").Value = 6

# --- Row 9: update values, add a new F9 cell, change G9 ------------------

$ws.Range("B9").Value = -0.1929832571191242
$ws.Range("C9").Value = 0.4511597977240172
$ws.Range("D9").Value = 0.2784649223588765
$ws.Range("E9").Value = 0.5276977566361984
$ws.Range("F9").Value = 0.6015260403721985
$ws.Range("G9").Value = 3

# --- Row 10: brand-new row (Q8) -------------------------------------------

$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = -0.2080565315694543
$ws.Range("C10").Value = 0.2080565315694543
$ws.Range("D10").Value = 0.04328752032871131
$ws.Range("E10").Value = 0.2080565315694543
$ws.Range("G10").Value = 1
